$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some "Price" column values in this sheet are digit-grouped numbers stored
# as plain text (e.g. "69.196.81", "1.00") rather than real numbers. Cells
# whose new text is a single plain decimal (e.g. "1.00", "7.43") would
# otherwise be auto-converted to a number by Excel on assignment, so those
# specific cells are pre-formatted as Text, written, then the temporary
# formatting is cleared again so the cell keeps the default (unstyled) look.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "69.196.81"
$ws.Range("E2").Value = "  +0.60%  "

# Row 3
$ws.Range("D3").Value = "3.778.12"
$ws.Range("E3").Value = "  +0.64%  "

# Row 4
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.04%  "

# Row 5
$ws.Range("D5").Value = "601.96"
$ws.Range("E5").Value = "  +0.03%  "

# Row 6
$ws.Range("D6").Value = "166.07"
$ws.Range("E6").Value = "  -1.84%  "

# Row 7
$ws.Range("D7").Value = "3.782.27"
$ws.Range("E7").Value = "  +0.83%  "

# Row 8
$ws.Range("E8").Value = "  -0.08%  "

# Row 9
$ws.Range("E9").Value = "  +0.48%  "

# Row 10
$ws.Range("E10").Value = "  +4.26%  "

# Row 11
$ws.Range("E11").Value = "  +0.29%  "

# Row 12
$ws.Range("D12").Value = "0.461"
$ws.Range("E12").Value = "  -0.26%  "

# Row 13
$ws.Range("D13").Value = "37.73"
$ws.Range("E13").Value = "  -1.53%  "

# Row 14
$ws.Range("E14").Value = "  +0.17%  "

# Row 15
$ws.Range("D15").Value = "4.405.63"
$ws.Range("E15").Value = "  +0.51%  "

# Row 16
$ws.Range("D16").Value = "3.775.13"
$ws.Range("E16").Value = "  +0.55%  "

# Row 17
$ws.Range("D17").Value = "69.300.51"
$ws.Range("E17").Value = "  +0.76%  "

# Row 18
$ws.Range("D18").Value = "7.43"
$ws.Range("E18").Value = "  +2.02%  "

# Row 19
$ws.Range("D19").Value = "17.66"
$ws.Range("E19").Value = "  +3.33%  "

# Row 20
$ws.Range("E20").Value = "  -0.98%  "

# Row 21
$ws.Range("D21").Value = "11.39"
$ws.Range("E21").Value = "  +4.81%  "

# Row 22
$ws.Range("D22").Value = "493.69"
$ws.Range("E22").Value = "  -0.40%  "

# Row 23
$ws.Range("D23").Value = "0.727"
$ws.Range("E23").Value = "  -0.38%  "

# Row 24
$ws.Range("E24").Value = "  -1.57%  "

# Row 25
$ws.Range("D25").Value = "84.94"
$ws.Range("E25").Value = "  -0.44%  "

# Row 26
$ws.Range("E26").Value = "  -2.51%  "

# Row 27
$ws.Range("D27").Value = "12.30"
$ws.Range("E27").Value = "  -0.48%  "

# Row 28
$ws.Range("D28").Value = "10.12"
$ws.Range("E28").Value = "  -1.68%  "

# Row 29
$ws.Range("E29").Value = "  +0.08%  "

# Row 30
$ws.Range("E30").Value = "  -0.51%  "

# Row 31
$ws.Range("E31").Value = "  +2.66%  "

# Row 32
$ws.Range("E32").Value = "  -3.69%  "

# Row 33
$ws.Range("D33").Value = "31.96"
$ws.Range("E33").Value = "  +0.14%  "

# Row 34
$ws.Range("D34").Value = "3.915.39"
$ws.Range("E34").Value = "  +0.38%  "

# Row 35
$ws.Range("D35").Value = "3.722.80"
$ws.Range("E35").Value = "  +0.91%  "

# Row 36
$ws.Range("E36").Value = "  -0.75%  "

# Row 37
$ws.Range("B37").Value = "Mantle"
$ws.Range("C37").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D37").Value = "1.02"
$ws.Range("E37").Value = "  +0.16%  "

# Row 38
$ws.Range("B38").Value = "Filecoin"
$ws.Range("C38").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D38").Value = "5.96"
$ws.Range("E38").Value = "  +1.75%  "

# Row 39
$ws.Range("D39").Value = "0.139"
$ws.Range("E39").Value = "  +4.43%  "

# Row 40
$ws.Range("E40").Value = "  +0.05%  "

# Row 41
$ws.Range("E41").Value = "  +5.92%  "

# Row 42
$ws.Range("E42").Value = "  +0.54%  "

# Row 43
$ws.Range("D43").Value = "48.56"
$ws.Range("E43").Value = "  -0.64%  "

# Row 44
$ws.Range("D44").Value = "426.49"
$ws.Range("E44").Value = "  -2.93%  "

# Row 45
$ws.Range("E45").Value = "  +0.63%  "

# Row 46
$ws.Range("D46").Value = "8.45"
$ws.Range("E46").Value = "  -0.33%  "

# Row 47
$ws.Range("E47").Value = "  -0.01%  "

# Row 48
$ws.Range("B48").Value = "Arweave"
$ws.Range("C48").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D48").Value = "40.12"
$ws.Range("E48").Value = "  -1.13%  "

# Row 49
$ws.Range("B49").Value = "Monero"
$ws.Range("C49").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D49").Value = "142.15"
$ws.Range("E49").Value = "  +0.40%  "

# Row 50
$ws.Range("D50").Value = "2.812.87"
$ws.Range("E50").Value = "  +0.86%  "

# Row 51
$ws.Range("D51").Value = "1.30"
$ws.Range("E51").Value = "  +8.96%  "

# Clear the temporary Text formatting residue.
$ws.Range("D4").ClearFormats()
$ws.Range("D5").ClearFormats()
$ws.Range("D6").ClearFormats()
$ws.Range("D12").ClearFormats()
$ws.Range("D13").ClearFormats()
$ws.Range("D18").ClearFormats()
$ws.Range("D19").ClearFormats()
$ws.Range("D21").ClearFormats()
$ws.Range("D22").ClearFormats()
$ws.Range("D23").ClearFormats()
$ws.Range("D25").ClearFormats()
$ws.Range("D27").ClearFormats()
$ws.Range("D28").ClearFormats()
$ws.Range("D33").ClearFormats()
$ws.Range("D37").ClearFormats()
$ws.Range("D38").ClearFormats()
$ws.Range("D39").ClearFormats()
$ws.Range("D43").ClearFormats()
$ws.Range("D44").ClearFormats()
$ws.Range("D46").ClearFormats()
$ws.Range("D48").ClearFormats()
$ws.Range("D49").ClearFormats()
$ws.Range("D51").ClearFormats()
